# Parallel testing support: add InvalidLoginData and ParallelLoginData sheets,
# and adjust the existing LoginData sheet's selection.

$wb = $excel.ActiveWorkbook
$loginData = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# InvalidLoginData (new sheet, placed right after LoginData)
# ---------------------------------------------------------------------------
$invalidLoginData = $wb.Worksheets.Add($null, $loginData)
$invalidLoginData.Name = "InvalidLoginData"

$invalidLoginData.Range("A1").Value = "Email"
$invalidLoginData.Range("B1").Value = "Password"
$invalidLoginData.Range("A1:B1").Interior.Color = 65535

$invalidLoginData.Range("A2").Value = "wronguser@example.com"
$invalidLoginData.Range("B2").Value = "Test123456"

$invalidLoginData.Range("A3").Value = "testuser@example.com"
$invalidLoginData.Range("B3").Value = "wrongpass"

$invalidLoginData.Range("B4").Value = "Test123456"

$invalidLoginData.Range("A5").Value = "testuser@example.com"

$invalidLoginData.Range("A7").Value = "test+user@example.com"
$invalidLoginData.Range("B7").Value = "Test@123"

$invalidLoginData.Columns.Item(1).ColumnWidth = 21.43
$invalidLoginData.Columns.Item(2).ColumnWidth = 9.5

$invalidLoginData.Range("A1:C3").Select()

# ---------------------------------------------------------------------------
# ParallelLoginData (new sheet, placed right after InvalidLoginData)
# ---------------------------------------------------------------------------
$parallelLoginData = $wb.Worksheets.Add($null, $invalidLoginData)
$parallelLoginData.Name = "ParallelLoginData"

$parallelLoginData.Range("A1").Value = "Email"
$parallelLoginData.Range("B1").Value = "Password"
$parallelLoginData.Range("A1:B1").Interior.Color = 65535

$parallelLoginData.Range("A2").Value = "wronguser@example.com"
$parallelLoginData.Range("B2").Value = "Test123456"

$parallelLoginData.Range("A3").Value = "testuser@example.com"
$parallelLoginData.Range("B3").Value = "wrongpass"

$parallelLoginData.Range("A4").Value = "wronguserexample.com"
$parallelLoginData.Range("B4").Value = "Test123456"

$parallelLoginData.Range("A5").Value = "testuser@example.com"

$parallelLoginData.Columns.Item(1).ColumnWidth = 21.43
$parallelLoginData.Columns.Item(2).ColumnWidth = 9.5

$parallelLoginData.Range("D13").Select()

# ---------------------------------------------------------------------------
# Update the original LoginData sheet's selection (no longer the active tab)
# ---------------------------------------------------------------------------
$loginData.Range("B1:C8").Select()

# ParallelLoginData ends up being the active / last-selected sheet.
$parallelLoginData.Activate()
$parallelLoginData.Range("D13").Select()
